$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 45949.773
$ws.Range("I135").Value = 514.4
$ws.Range("K135").Value = 4629.599999999999
$ws.Range("M135").Value = -2094.599999999999


# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 23086.5
$ws.Range("J75").Value = 23086.5
$ws.Range("L75").Value = 23086.5
$ws.Range("N75").Value = -24834.5

$ws.Range("H78").Value = 23086.5
$ws.Range("J78").Value = 23086.5
$ws.Range("L78").Value = 69259.5
$ws.Range("N78").Value = -77995.5

$ws.Range("H97").Value = 771
$ws.Range("I97").Value = 822.4167
$ws.Range("J97").Value = 462.5
$ws.Range("K97").Value = 822.4167
$ws.Range("L97").Value = 462.5
$ws.Range("M97").Value = -326.4167
$ws.Range("N97").Value = -1454.5

$ws.Range("H101").Value = 40600.332
$ws.Range("J101").Value = 40600.332
$ws.Range("L101").Value = 40600.332
$ws.Range("N101").Value = -47090.332

$ws.Range("H107").Value = 37614
$ws.Range("J107").Value = 37614
$ws.Range("L107").Value = 37614
$ws.Range("N107").Value = -45294

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 2231.0286
$ws.Range("I122").Value = 2169.926
$ws.Range("K122").Value = 6509.778
$ws.Range("M122").Value = -4059.778


# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 27480.75
$ws.Range("J100").Value = 27480.75
$ws.Range("L100").Value = 27480.75
$ws.Range("N100").Value = -29644.75

$ws.Range("H105").Value = 1313.1428
$ws.Range("I105").Value = 1313.1428
$ws.Range("K105").Value = 1313.1428
$ws.Range("M105").Value = 433.8571999999999


# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10758.228
$ws.Range("I31").Value = 3146.5
$ws.Range("J31").Value = 18369.955
$ws.Range("K31").Value = 3146.5
$ws.Range("L31").Value = 18369.955
$ws.Range("M31").Value = -2851.5
$ws.Range("N31").Value = -18959.955

$ws.Range("H34").Value = 10758.228
$ws.Range("I34").Value = 3146.5
$ws.Range("J34").Value = 18369.955
$ws.Range("K34").Value = 3146.5
$ws.Range("L34").Value = 18369.955
$ws.Range("M34").Value = -2944.5
$ws.Range("N34").Value = -18773.955

$ws.Range("H58").Value = 1884.1852
$ws.Range("I58").Value = 1850.2106
$ws.Range("J58").Value = 1964.875
$ws.Range("K58").Value = 1850.2106
$ws.Range("L58").Value = 1964.875
$ws.Range("M58").Value = -1647.2106
$ws.Range("N58").Value = -2370.875

$ws.Range("H62").Value = 4994.625
$ws.Range("I62").Value = 5598.2
$ws.Range("J62").Value = 3988.6667
$ws.Range("K62").Value = 5598.2
$ws.Range("L62").Value = 3988.6667
$ws.Range("M62").Value = -4974.2
$ws.Range("N62").Value = -5236.6667

$ws.Range("H65").Value = 4994.625
$ws.Range("I65").Value = 5598.2
$ws.Range("J65").Value = 3988.6667
$ws.Range("K65").Value = 27991
$ws.Range("L65").Value = 19943.3335
$ws.Range("M65").Value = -24871
$ws.Range("N65").Value = -26183.3335

$ws.Range("H92").Value = 27266.666
$ws.Range("J92").Value = 27266.666
$ws.Range("L92").Value = 27266.666
$ws.Range("N92").Value = -32258.666

$ws.Range("H96").Value = 11229.75
$ws.Range("J96").Value = 11229.75
$ws.Range("L96").Value = 11229.75
$ws.Range("N96").Value = -16721.75

$ws.Range("H132").Value = 3017.4285
$ws.Range("I132").Value = 3022.182
$ws.Range("K132").Value = 9066.545999999998
$ws.Range("M132").Value = -6536.545999999998

$ws.Range("H134").Value = 3495.3447
$ws.Range("I134").Value = 3469.625
$ws.Range("K134").Value = 10408.875
$ws.Range("M134").Value = -7873.875

$ws.Range("H136").Value = 1884.1852
$ws.Range("I136").Value = 1850.2106
$ws.Range("J136").Value = 1964.875
$ws.Range("K136").Value = 5550.6318
$ws.Range("L136").Value = 5894.625
$ws.Range("M136").Value = -3000.6318
$ws.Range("N136").Value = -10994.625


# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 486.86667
$ws.Range("I5").Value = 497.75
$ws.Range("J5").Value = 465.1
$ws.Range("K5").Value = 1493.25
$ws.Range("L5").Value = 1395.3
$ws.Range("M5").Value = -1381.25
$ws.Range("N5").Value = -1619.3

$ws.Range("H56").Value = 10237
$ws.Range("I56").Value = 10237
$ws.Range("K56").Value = 10237
$ws.Range("M56").Value = -9707

$ws.Range("H113").Value = 760.25806
$ws.Range("I113").Value = 703.1818
$ws.Range("J113").Value = 791.65
$ws.Range("K113").Value = 2109.5454
$ws.Range("L113").Value = 2374.95
$ws.Range("M113").Value = 60.45460000000003
$ws.Range("N113").Value = -6714.95

$ws.Range("H132").Value = 1718.0869
$ws.Range("I132").Value = 1176.7
$ws.Range("J132").Value = 2134.5386
$ws.Range("K132").Value = 10590.3
$ws.Range("L132").Value = 19210.8474
$ws.Range("M132").Value = -8060.300000000001
$ws.Range("N132").Value = -24270.8474

$ws.Range("H135").Value = 486.86667
$ws.Range("I135").Value = 497.75
$ws.Range("J135").Value = 465.1
$ws.Range("K135").Value = 4479.75
$ws.Range("L135").Value = 4185.900000000001
$ws.Range("M135").Value = -1944.75
$ws.Range("N135").Value = -9255.900000000001


# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 6833.3335
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 6833.3335
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 6833.3335
$ws.Range("N22").Value = -7891.3335
$ws.Range("M22").ClearContents()

$ws.Range("H97").Value = 26936.96
$ws.Range("I97").Value = 26936.96
$ws.Range("K97").Value = 26936.96
$ws.Range("M97").Value = -26440.96

$ws.Range("H122").Value = 1955.5333
$ws.Range("I122").Value = 1682.6
$ws.Range("J122").Value = 2501.4
$ws.Range("K122").Value = 5047.799999999999
$ws.Range("L122").Value = 7504.200000000001
$ws.Range("M122").Value = -2597.799999999999
$ws.Range("N122").Value = -12404.2


# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1354.5555
$ws.Range("J22").Value = 2666.3333
$ws.Range("L22").Value = 2666.3333
$ws.Range("N22").Value = -3256.3333

$ws.Range("H27").Value = 1354.5555
$ws.Range("J27").Value = 2666.3333
$ws.Range("L27").Value = 2666.3333
$ws.Range("N27").Value = -2880.3333

$ws.Range("H61").Value = 40582.223
$ws.Range("I61").Value = 3843.1428
$ws.Range("K61").Value = 3843.1428
$ws.Range("M61").Value = -3641.1428

$ws.Range("H68").Value = 1565.5
$ws.Range("I68").Value = 1511.375
$ws.Range("J68").Value = 1673.75
$ws.Range("K68").Value = 1511.375
$ws.Range("L68").Value = 1673.75
$ws.Range("M68").Value = -762.375
$ws.Range("N68").Value = -3171.75

$ws.Range("H71").Value = 1565.5
$ws.Range("I71").Value = 1511.375
$ws.Range("J71").Value = 1673.75
$ws.Range("K71").Value = 7556.875
$ws.Range("L71").Value = 8368.75
$ws.Range("M71").Value = -3812.875
$ws.Range("N71").Value = -15856.75

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H113").Value = 40582.223
$ws.Range("I113").Value = 3843.1428
$ws.Range("K113").Value = 3843.1428
$ws.Range("M113").Value = -1673.1428

$ws.Range("H132").Value = 3694.8696
$ws.Range("I132").Value = 3033.5557
$ws.Range("J132").Value = 4120
$ws.Range("K132").Value = 9100.667099999999
$ws.Range("L132").Value = 12360
$ws.Range("M132").Value = -6570.667099999999
$ws.Range("N132").Value = -17420


# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13668.846
$ws.Range("I62").Value = 14674
$ws.Range("K62").Value = 14674
$ws.Range("M62").Value = -14050

$ws.Range("H65").Value = 13668.846
$ws.Range("I65").Value = 14674
$ws.Range("K65").Value = 73370
$ws.Range("M65").Value = -70250

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H81").Value = 7195
$ws.Range("J81").Value = 13465.667
$ws.Range("L81").Value = 26931.334
$ws.Range("N81").Value = -29053.334

$ws.Range("H84").Value = 7195
$ws.Range("J84").Value = 13465.667
$ws.Range("L84").Value = 134656.67
$ws.Range("N84").Value = -145264.67

